# RN-1289: update the entity associated with a survey resubmission
# The "Clinic Code" / "Clinic Name" column headers used for the resubmission
# target row are renamed to "Entity Code" / "Entity Name" on the
# Test_Clinic_Data sheet (row 2 = D2, row 3 = D3).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test_Clinic_Data")

$ws.Range("D3").Value = "Entity Name"
$ws.Range("D2").Value = "Entity Code"

# Match the author's final active-cell selection on this sheet.
$ws.Range("D3").Select()
